# "Removed Hydric entries in the table"
#
# The VariableRankings sheet listed two "hydric" response variables
# (PCT_HYDRIC / AREA_HYDRIC) in rows 81-82. Their data is cleared out
# (the rows are kept, but all the cell content in them is removed),
# which also drops the shared strings that are no longer referenced and
# updates the dependent COUNTIF() formulas on the Metadata sheet.
#
# A note documenting the change is appended to the Notes sheet, and the
# Notes sheet becomes the active/selected sheet (matching the workbook's
# saved view state after the edit).

$wb = $excel.ActiveWorkbook

# --- VariableRankings: clear out the two hydric rows ---
$rankings = $wb.Worksheets.Item("VariableRankings")
$rankings.Range("A81:M82").ClearContents()

# Minor row-height tweak that was present in the saved workbook.
$rankings.Rows.Item(80).RowHeight = 13.5

# Leave the selection where it ended up after removing those rows.
$rankings.Range("C84").Select()

# --- Notes: log the change and make it the active sheet ---
$notes = $wb.Worksheets.Item("Notes")
$notes.Range("B5").Value = "Removed Hydric Layers"
$notes.Activate()
